$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide "Prerequisites" (sldId 260, title shape creationId
#    {C566ADFD-99F0-C9A1-3DC1-57CAEB3749CA}) - nudge the Title placeholder's
#    position by giving it an explicit xfrm override. Top moves from the
#    inherited 609600 EMU up to 589280 EMU (46.4pt instead of 48pt); left /
#    width / height stay at the values already inherited from the layout.
# ---------------------------------------------------------------------------
$targetSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ($slide.SlideID -eq 260) {
        $targetSlide = $slide
        break
    }
}
if ($targetSlide -eq $null) {
    # Fallback: find it by its title text if SlideID lookup didn't resolve.
    for ($i = 1; $i -le $p.Slides.Count; $i++) {
        $slide = $p.Slides.Item($i)
        if ($slide.Shapes.HasTitle) {
            if ($slide.Shapes.Title.TextFrame.TextRange.Text -eq "Prerequisites") {
                $targetSlide = $slide
                break
            }
        }
    }
}
if ($targetSlide -ne $null -and $targetSlide.Shapes.HasTitle) {
    $targetSlide.Shapes.Title.Top = 46.4
}

# ---------------------------------------------------------------------------
# 2) Refresh the "update automatically" date field cached on every slide
#    layout and the slide master (last-save date stamp: 6/5/2023 ->
#    8/28/2023).
# ---------------------------------------------------------------------------
$newDate = "8/28/2023"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

Write-Output "Edit complete"
